# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit-tracking workbook
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 29.666666
$ws.Range("I39").Value = 38.25
$ws.Range("J39").Value = 12.5
$ws.Range("K39").Value = 114.75
$ws.Range("L39").Value = 37.5
$ws.Range("M39").Value = 181.25
$ws.Range("N39").Value = -629.5
$ws.Range("H40").Value = 4002
$ws.Range("J40").Value = 4002
$ws.Range("L40").Value = 4002
$ws.Range("N40").Value = -4352
$ws.Range("H51").Value = 7999.6665
$ws.Range("I51").Value = 7999
$ws.Range("K51").Value = 7999
$ws.Range("M51").Value = -7515
$ws.Range("H58").Value = 2535.6
$ws.Range("I58").Value = 467
$ws.Range("J58").Value = 3422.1428
$ws.Range("K58").Value = 1401
$ws.Range("L58").Value = 10266.4284
$ws.Range("M58").Value = -1251
$ws.Range("N58").Value = -10566.4284
$ws.Range("H76").Value = 3903
$ws.Range("I76").Value = 3903
$ws.Range("K76").Value = 3903
$ws.Range("M76").Value = -3588
$ws.Range("H79").Value = 3903
$ws.Range("I79").Value = 3903
$ws.Range("K79").Value = 3903
$ws.Range("M79").Value = -2811
$ws.Range("H137").Value = 530331.3
$ws.Range("I137").Value = 716235.1
$ws.Range("K137").Value = 2148705.3
$ws.Range("M137").Value = -2146155.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 998
$ws.Range("J12").Value = 998
$ws.Range("L12").Value = 998
$ws.Range("N12").Value = -1344
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 111111
$ws.Range("J112").Value = 111111
$ws.Range("L112").Value = 111111
$ws.Range("N112").Value = -114065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 799
$ws.Range("I22").Value = 799
$ws.Range("K22").Value = 799
$ws.Range("M22").Value = -449
$ws.Range("H43").Value = 16856.6
$ws.Range("J43").Value = 16856.6
$ws.Range("L43").Value = 16856.6
$ws.Range("N43").Value = -17224.6
$ws.Range("H62").Value = 4232.8335
$ws.Range("I62").Value = 3874.25
$ws.Range("K62").Value = 3874.25
$ws.Range("M62").Value = -3250.25
$ws.Range("H65").Value = 4232.8335
$ws.Range("I65").Value = 3874.25
$ws.Range("K65").Value = 19371.25
$ws.Range("M65").Value = -16251.25
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H99").Value = 5803.6
$ws.Range("I99").Value = 5803.6
$ws.Range("K99").Value = 5803.6
$ws.Range("M99").Value = -4305.6
$ws.Range("H101").Value = 16856.6
$ws.Range("J101").Value = 16856.6
$ws.Range("L101").Value = 16856.6
$ws.Range("N101").Value = -23346.6
$ws.Range("H107").Value = 556.53845
$ws.Range("I107").Value = 576.2727
$ws.Range("J107").Value = 448
$ws.Range("K107").Value = 576.2727
$ws.Range("L107").Value = 448
$ws.Range("M107").Value = 1343.7273
$ws.Range("N107").Value = -4288
$ws.Range("H126").Value = 5803.6
$ws.Range("I126").Value = 5803.6
$ws.Range("K126").Value = 17410.8
$ws.Range("M126").Value = -14940.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 207.6
$ws.Range("I2").Value = 267.26315
$ws.Range("J2").Value = 18.666666
$ws.Range("K2").Value = 1603.5789
$ws.Range("L2").Value = 111.999996
$ws.Range("M2").Value = -1490.5789
$ws.Range("N2").Value = -337.999996
$ws.Range("H4").Value = 100011000
$ws.Range("J4").Value = 100011000
$ws.Range("L4").Value = 300033000
$ws.Range("N4").Value = -300033224
$ws.Range("H17").Value = 716.25
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -431
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H29").Value = 31
$ws.Range("I29").Value = 60
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 180
$ws.Range("L29").Value = 6
$ws.Range("M29").Value = 97
$ws.Range("N29").Value = -560
$ws.Range("H34").Value = 350.5
$ws.Range("I34").Value = 301
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 903
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -819
$ws.Range("N34").Value = -1368
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H50").Value = 106.666664
$ws.Range("I50").Value = 106.666664
$ws.Range("K50").Value = 319.999992
$ws.Range("M50").Value = 161.000008
$ws.Range("H53").Value = 106.666664
$ws.Range("I53").Value = 106.666664
$ws.Range("K53").Value = 319.999992
$ws.Range("M53").Value = 161.000008
$ws.Range("H55").Value = 1000
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("H131").Value = 2657.182
$ws.Range("I131").Value = 1000
$ws.Range("K131").Value = 3000
$ws.Range("M131").Value = 2040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H23").Value = 250
$ws.Range("I23").Value = 250
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 250
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -27
$ws.Range("N23").ClearContents()
$ws.Range("H102").Value = 1277.625
$ws.Range("I102").Value = 1277.625
$ws.Range("K102").Value = 1277.625
$ws.Range("M102").Value = 344.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 729.4
$ws.Range("J46").Value = 729.4
$ws.Range("L46").Value = 729.4
$ws.Range("N46").Value = -1105.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 169333.33
$ws.Range("I7").Value = 169333.33
$ws.Range("K7").Value = 169333.33
$ws.Range("M7").Value = -169220.33
$ws.Range("H62").Value = 56250.5
$ws.Range("I62").Value = 8334
$ws.Range("K62").Value = 8334
$ws.Range("M62").Value = -7710
$ws.Range("H65").Value = 56250.5
$ws.Range("I65").Value = 8334
$ws.Range("K65").Value = 41670
$ws.Range("M65").Value = -38550
$ws.Range("H68").Value = 38271
$ws.Range("J68").Value = 38271
$ws.Range("L68").Value = 38271
$ws.Range("N68").Value = -39893
$ws.Range("H71").Value = 38271
$ws.Range("J71").Value = 38271
$ws.Range("L71").Value = 114813
$ws.Range("N71").Value = -122925
$ws.Range("H81").Value = 1225
$ws.Range("I81").Value = 1199.3334
$ws.Range("J81").Value = 1302
$ws.Range("K81").Value = 2398.6668
$ws.Range("L81").Value = 2604
$ws.Range("M81").Value = -1337.6668
$ws.Range("N81").Value = -4726
$ws.Range("H84").Value = 1225
$ws.Range("I84").Value = 1199.3334
$ws.Range("J84").Value = 1302
$ws.Range("K84").Value = 11993.334
$ws.Range("L84").Value = 13020
$ws.Range("M84").Value = -6689.333999999999
$ws.Range("N84").Value = -23628
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H136").Value = 7856.857
$ws.Range("I136").Value = 3249.75
$ws.Range("K136").Value = 9749.25
$ws.Range("M136").Value = -7199.25
